$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose target price text is a "clean" number (e.g. "258.07") need the
# number format forced to Text first, otherwise Excel auto-converts the typed
# text into a numeric value (losing significant trailing zeros, e.g. "12.60").
$textForceCells = @("D5","D6","D8","D13","D19","D20","D22","D23","D25","D26","D28","D33","D34","D35","D36","D37","D39","D40","D42","D43","D45","D48","D49","D50")
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$updates = @(
    @{ Cell = "D2"; Value = "98.479.02" },
    @{ Cell = "E2"; Value = "  -0.11%  " },
    @{ Cell = "D3"; Value = "3.366.93" },
    @{ Cell = "E4"; Value = "  +0.01%  " },
    @{ Cell = "D5"; Value = "258.07" },
    @{ Cell = "D6"; Value = "666.32" },
    @{ Cell = "E6"; Value = "  +5.91%  " },
    @{ Cell = "E7"; Value = "  +10.23%  " },
    @{ Cell = "D8"; Value = "0.465" },
    @{ Cell = "E8"; Value = "  +19.41%  " },
    @{ Cell = "E9"; Value = "  +25.25%  " },
    @{ Cell = "E10"; Value = "  -0.01%  " },
    @{ Cell = "D11"; Value = "3.364.40" },
    @{ Cell = "E11"; Value = "  +0.21%  " },
    @{ Cell = "E12"; Value = "  +6.40%  " },
    @{ Cell = "D13"; Value = "42.03" },
    @{ Cell = "E13"; Value = "  +15.09%  " },
    @{ Cell = "E14"; Value = "  +9.30%  " },
    @{ Cell = "D15"; Value = "96.771.12" },
    @{ Cell = "E15"; Value = "  -1.66%  " },
    @{ Cell = "D16"; Value = "3.999.14" },
    @{ Cell = "E16"; Value = "  +0.27%  " },
    @{ Cell = "E17"; Value = "  +3.12%  " },
    @{ Cell = "D18"; Value = "3.361.43" },
    @{ Cell = "E18"; Value = "  -0.27%  " },
    @{ Cell = "D19"; Value = "7.61" },
    @{ Cell = "E19"; Value = "  +25.08%  " },
    @{ Cell = "D20"; Value = "16.82" },
    @{ Cell = "E20"; Value = "  +10.59%  " },
    @{ Cell = "E21"; Value = "  +1.02%  " },
    @{ Cell = "D22"; Value = "530.46" },
    @{ Cell = "E22"; Value = "  +8.48%  " },
    @{ Cell = "D23"; Value = "10.54" },
    @{ Cell = "E23"; Value = "  +12.19%  " },
    @{ Cell = "E24"; Value = "  +4.07%  " },
    @{ Cell = "D25"; Value = "0.438" },
    @{ Cell = "E25"; Value = "  +55.06%  " },
    @{ Cell = "D26"; Value = "102.74" },
    @{ Cell = "E26"; Value = "  +15.25%  " },
    @{ Cell = "E27"; Value = "  +10.79%  " },
    @{ Cell = "D28"; Value = "12.60" },
    @{ Cell = "E28"; Value = "  +5.72%  " },
    @{ Cell = "D29"; Value = "3.545.96" },
    @{ Cell = "E29"; Value = "  -0.08%  " },
    @{ Cell = "E30"; Value = "  +10.73%  " },
    @{ Cell = "E31"; Value = "  -0.41%  " },
    @{ Cell = "E32"; Value = "  +14.12%  " },
    @{ Cell = "D33"; Value = "0.191" },
    @{ Cell = "E33"; Value = "  +0.05%  " },
    @{ Cell = "D34"; Value = "1.00" },
    @{ Cell = "E34"; Value = "  +0.10%  " },
    @{ Cell = "D35"; Value = "29.55" },
    @{ Cell = "E35"; Value = "  +5.02%  " },
    @{ Cell = "D36"; Value = "0.548" },
    @{ Cell = "E36"; Value = "  +19.16%  " },
    @{ Cell = "D37"; Value = "7.88" },
    @{ Cell = "E37"; Value = "  +8.26%  " },
    @{ Cell = "E38"; Value = "  +8.52%  " },
    @{ Cell = "D39"; Value = "0.158" },
    @{ Cell = "E39"; Value = "  +5.53%  " },
    @{ Cell = "D40"; Value = "530.39" },
    @{ Cell = "E40"; Value = "  +6.16%  " },
    @{ Cell = "E41"; Value = "  +6.74%  " },
    @{ Cell = "D42"; Value = "24.70" },
    @{ Cell = "E42"; Value = "  -0.55%  " },
    @{ Cell = "D43"; Value = "0.0438" },
    @{ Cell = "E43"; Value = "  +34.38%  " },
    @{ Cell = "E44"; Value = "  +1.24%  " },
    @{ Cell = "D45"; Value = "3.44" },
    @{ Cell = "E45"; Value = "  +4.72%  " },
    @{ Cell = "E46"; Value = "  +5.18%  " },
    @{ Cell = "E47"; Value = "  -0.01%  " },
    @{ Cell = "D48"; Value = "2.08" },
    @{ Cell = "E48"; Value = "  +7.33%  " },
    @{ Cell = "D49"; Value = "7.95" },
    @{ Cell = "E49"; Value = "  +18.70%  " },
    @{ Cell = "D50"; Value = "5.12" },
    @{ Cell = "E50"; Value = "  +10.78%  " },
    @{ Cell = "E51"; Value = "  +11.93%  " }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
